# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This updates the "K" column (column G) values for each game-log row with
# the recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of worksheet row number -> new K (strikeouts) value.
$kValues = @{
    2  = 4
    3  = 1
    4  = 1
    5  = 3
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 2
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 2
    29 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 3
    36 = 1
    37 = 3
    38 = 2
    39 = 1
    40 = 3
    41 = 2
    42 = 1
    43 = 1
    44 = 1
    45 = 0
    46 = 2
    47 = 2
    48 = 0
    49 = 0
    50 = 2
    51 = 1
    52 = 0
    53 = 1
    54 = 1
    55 = 1
    57 = 0
    58 = 1
    59 = 2
    60 = 1
    61 = 1
    62 = 0
    63 = 3
    64 = 2
    65 = 0
    66 = 2
    67 = 1
    68 = 0
    69 = 1
    70 = 3
    71 = 2
    72 = 2
    73 = 2
    74 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
